# Apply the edit described by the diff:
# - Extend the table from columns A:O to A:Q (add columns P and Q).
# - Row 1 (header row): add P1=14, Q1=15 with the same style as O1.
# - Rows 2-25 (data rows): flip I<->K and M<->O values (I:1->2, K:2->1,
#   M:1->2, O:2->1), and add new columns P=2, Q=2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 and Q1, matching O1's style ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P (new) = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q (new) = 2
}
